# Auto-generated Excel COM-interop script applying the Leviathan_Profits market-data refresh
# (scheduled runner update): updates currentAveragePrice / Leve price / profit columns
# per worksheet (crafting-class tabs) to the latest fetched market values.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 222.54546
$ws.Range("I28").Value = 232.19048
$ws.Range("J28").Value = 20
$ws.Range("K28").Value = 232.19048
$ws.Range("L28").Value = 20
$ws.Range("M28").Value = 252.80952
$ws.Range("N28").Value = -990
$ws.Range("H57").Value = 50435.57
$ws.Range("J57").Value = 50435.57
$ws.Range("L57").Value = 151306.71
$ws.Range("N57").Value = -152304.71
$ws.Range("H61").Value = 495
$ws.Range("I61").Value = 495
$ws.Range("K61").Value = 1485
$ws.Range("M61").Value = -1313
$ws.Range("H97").Value = 1639.8572
$ws.Range("J97").Value = 2020.4
$ws.Range("L97").Value = 6061.200000000001
$ws.Range("N97").Value = -7053.200000000001
$ws.Range("H99").Value = 66667084
$ws.Range("I99").Value = 482.55554
$ws.Range("J99").Value = 166666990
$ws.Range("K99").Value = 1447.66662
$ws.Range("L99").Value = 500000970
$ws.Range("M99").Value = 50.33338000000003
$ws.Range("N99").Value = -500003966
$ws.Range("H101").Value = 958.8889
$ws.Range("I101").Value = 970.73334
$ws.Range("K101").Value = 2912.20002
$ws.Range("M101").Value = -1290.20002
$ws.Range("H127").Value = 937
$ws.Range("I127").Value = 926.5
$ws.Range("J127").Value = 1000
$ws.Range("K127").Value = 2779.5
$ws.Range("L127").Value = 3000
$ws.Range("M127").Value = 2180.5
$ws.Range("N127").Value = -12920
$ws.Range("H138").Value = 2317.5833
$ws.Range("I138").Value = 1954.5625
$ws.Range("J138").Value = 2608
$ws.Range("K138").Value = 5863.6875
$ws.Range("L138").Value = 7824
$ws.Range("M138").Value = -723.6875
$ws.Range("N138").Value = -18104
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H138").Value = 100129.336
$ws.Range("I138").Value = 60390
$ws.Range("J138").Value = 119999
$ws.Range("K138").Value = 60390
$ws.Range("L138").Value = 119999
$ws.Range("M138").Value = -55250
$ws.Range("N138").Value = -130279
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1369.8064
$ws.Range("I134").Value = 1142.8
$ws.Range("K134").Value = 3428.4
$ws.Range("M134").Value = -893.3999999999996
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10032.458
$ws.Range("I31").Value = 2829.0857
$ws.Range("J31").Value = 29426.154
$ws.Range("K31").Value = 2829.0857
$ws.Range("L31").Value = 29426.154
$ws.Range("M31").Value = -2534.0857
$ws.Range("N31").Value = -30016.154
$ws.Range("H34").Value = 10032.458
$ws.Range("I34").Value = 2829.0857
$ws.Range("J34").Value = 29426.154
$ws.Range("K34").Value = 2829.0857
$ws.Range("L34").Value = 29426.154
$ws.Range("M34").Value = -2627.0857
$ws.Range("N34").Value = -29830.154
$ws.Range("H107").Value = 1571.2693
$ws.Range("J107").Value = 1371.1428
$ws.Range("L107").Value = 1371.1428
$ws.Range("N107").Value = -5211.1428
$ws.Range("H134").Value = 3958.6667
$ws.Range("I134").Value = 2228.2
$ws.Range("K134").Value = 6684.599999999999
$ws.Range("M134").Value = -4149.599999999999
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1198.3889
$ws.Range("I5").Value = 795.7273
$ws.Range("K5").Value = 2387.1819
$ws.Range("M5").Value = -2275.1819
$ws.Range("H52").Value = 1422.5
$ws.Range("J52").Value = 1422.5
$ws.Range("L52").Value = 4267.5
$ws.Range("N52").Value = -4799.5
$ws.Range("H56").Value = 9332
$ws.Range("I56").Value = 9332
$ws.Range("K56").Value = 9332
$ws.Range("M56").Value = -8802
$ws.Range("H70").Value = 4249.8335
$ws.Range("I70").Value = 2501
$ws.Range("J70").Value = 5124.25
$ws.Range("K70").Value = 7503
$ws.Range("L70").Value = 15372.75
$ws.Range("M70").Value = -7188
$ws.Range("N70").Value = -16002.75
$ws.Range("H73").Value = 4249.8335
$ws.Range("I73").Value = 2501
$ws.Range("J73").Value = 5124.25
$ws.Range("K73").Value = 7503
$ws.Range("L73").Value = 15372.75
$ws.Range("M73").Value = -6411
$ws.Range("N73").Value = -17556.75
$ws.Range("H88").Value = 10690.5
$ws.Range("J88").Value = 10690.5
$ws.Range("L88").Value = 32071.5
$ws.Range("N88").Value = -32927.5
$ws.Range("H91").Value = 10690.5
$ws.Range("J91").Value = 10690.5
$ws.Range("L91").Value = 32071.5
$ws.Range("N91").Value = -35035.5
$ws.Range("H135").Value = 1198.3889
$ws.Range("I135").Value = 795.7273
$ws.Range("K135").Value = 7161.545700000001
$ws.Range("M135").Value = -4626.545700000001
$ws.Range("H137").Value = 2534.5
$ws.Range("I137").Value = 1824.2307
$ws.Range("J137").Value = 3853.5715
$ws.Range("K137").Value = 5472.6921
$ws.Range("L137").Value = 11560.7145
$ws.Range("M137").Value = -372.6921000000002
$ws.Range("N137").Value = -21760.7145
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("H97").Value = 41680.125
$ws.Range("I97").Value = 26643.68
$ws.Range("J97").Value = 95381.71000000001
$ws.Range("K97").Value = 26643.68
$ws.Range("L97").Value = 95381.71000000001
$ws.Range("M97").Value = -26147.68
$ws.Range("N97").Value = -96373.71000000001
$ws.Range("H107").Value = 38466670
$ws.Range("I107").Value = 295.4
$ws.Range("J107").Value = 62508156
$ws.Range("K107").Value = 295.4
$ws.Range("L107").Value = 62508156
$ws.Range("M107").Value = 1624.6
$ws.Range("N107").Value = -62511996
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 24911.895
$ws.Range("I46").Value = 55103
$ws.Range("J46").Value = 2954.7273
$ws.Range("K46").Value = 55103
$ws.Range("L46").Value = 2954.7273
$ws.Range("M46").Value = -54915
$ws.Range("N46").Value = -3330.7273
$ws.Range("H48").Value = 37495
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("H55").Value = 552
$ws.Range("J55").Value = 378.33334
$ws.Range("L55").Value = 378.33334
$ws.Range("N55").Value = -724.33334
$ws.Range("H61").Value = 72351.16
$ws.Range("I61").Value = 78718.69500000001
$ws.Range("K61").Value = 78718.69500000001
$ws.Range("M61").Value = -78516.69500000001
$ws.Range("H93").Value = 10245
$ws.Range("I93").Value = 1530.9714
$ws.Range("K93").Value = 1530.9714
$ws.Range("M93").Value = -282.9713999999999
$ws.Range("H113").Value = 72351.16
$ws.Range("I113").Value = 78718.69500000001
$ws.Range("K113").Value = 78718.69500000001
$ws.Range("M113").Value = -76548.69500000001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2678.475
$ws.Range("I136").Value = 2449.4
$ws.Range("J136").Value = 3365.7
$ws.Range("K136").Value = 7348.200000000001
$ws.Range("L136").Value = 10097.1
$ws.Range("M136").Value = -4798.200000000001
$ws.Range("N136").Value = -15197.1

# Two rows lost a now-inapplicable HQ/NQ profit column entirely (cell removed, not zeroed)
$wsGSM = $wb.Worksheets.Item("GSM")
$wsGSM.Range("M44").ClearContents()
$wsLTW = $wb.Worksheets.Item("LTW")
$wsLTW.Range("M48").ClearContents()

